# Generate Report for Handback
# Update timestamps / status text on the handback-status report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for both data rows
$wsOverview.Range("G2").Value = "2016-10-13 14:11:13"
$wsOverview.Range("G3").Value = "2016-10-13 14:11:13"

# zh-cn sheet
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-10-13 14:11:01"
$wsZhCn.Range("H3").Value = "2016-10-13 14:11:01"

# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-10-13 14:11:44"
$wsZhCn.Range("K3").Value = "2016-10-13 14:11:44"

# de-de sheet
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H2").Value = "2016-10-13 14:11:13"
$wsDeDe.Range("H3").Value = "2016-10-13 14:11:13"

# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-10-13 14:12:00"
$wsDeDe.Range("K3").Value = "2016-10-13 14:12:00"
